# feat: add 2022-Q4 data
#
# Before: sheet "总计" (totals) + sheet "2022-Q3" (fund holdings detail)
# After:  sheet "总计" (totals, now listing both Q4 and Q3 rows)
#       + sheet "2022-Q4" (fund holdings detail for the new quarter -- this
#         reuses the worksheet that used to be named "2022-Q3")
#       + sheet "2022-Q3" (fund holdings detail for Q3, a new worksheet
#         that is an exact copy of the data that used to live on the
#         "2022-Q3" tab before that tab got repurposed for Q4)

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Item(2)

# Step 1: duplicate the existing "2022-Q3" detail sheet to the end of the
# workbook *before* overwriting it with Q4 data, so the old Q3 detail data
# ends up preserved on its own tab. Rename the original first (so the name
# "2022-Q3" is free for the new copy to take).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$q3Sheet.Copy($null, $lastSheet)
$newQ3Sheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$q3Sheet.Name = "2022-Q4"
$newQ3Sheet.Name = "2022-Q3"

$q4Sheet = $q3Sheet

# Step 2: clear the (now-named) "2022-Q4" sheet and fill it with the new
# quarter's fund-holding detail rows, re-using the bold/bordered header
# style from the "总计" sheet (style index "2" in this workbook).
$q4Sheet.Cells.Clear()

# Match the page margins used on the "总计" sheet (0.75in/1in/0.5in)
# rather than the ones inherited from the old Q3 detail sheet.
$q4Sheet.PageSetup.LeftMargin = 54
$q4Sheet.PageSetup.RightMargin = 54
$q4Sheet.PageSetup.TopMargin = 72
$q4Sheet.PageSetup.BottomMargin = 72
$q4Sheet.PageSetup.HeaderMargin = 36
$q4Sheet.PageSetup.FooterMargin = 36

$styleSrc = $totalSheet.Cells.Item(1, 2)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $headers) {
    $styleSrc.Copy()
    $q4Sheet.Cells.Item(1, $col).PasteSpecial(-4122)
    $q4Sheet.Cells.Item(1, $col).Value2 = $h
    $col = $col + 1
}

$q4Rows = @(
    @("0", "014841", "东方阿尔法医疗健康混合A", "1.58", "89.34", "7.85", "0.1240", 2),
    @("1", "011471", "鹏华致远成长混合A",       "1.67", "60.84", "2.03", "0.0339", 2),
    @("2", "014842", "东方阿尔法医疗健康混合C", "0.43", "89.34", "7.85", "0.0338", 2),
    @("3", "011472", "鹏华致远成长混合C",       "0.06", "60.84", "2.03", "0.0012", 2)
)

$rowIdx = 2
foreach ($rowData in $q4Rows) {
    $styleSrc.Copy()
    $q4Sheet.Cells.Item($rowIdx, 1).PasteSpecial(-4122)
    $q4Sheet.Cells.Item($rowIdx, 1).Value2 = [int]$rowData[0]

    # Columns B-G are all stored as text (fund codes / names / numeric-
    # looking ratios kept as strings), so force a Text number format
    # before assigning, otherwise Excel would "helpfully" reinterpret
    # things like "014841" or "1.58" as numbers.
    $textRange = $q4Sheet.Range($q4Sheet.Cells.Item($rowIdx, 2), $q4Sheet.Cells.Item($rowIdx, 7))
    $textRange.NumberFormat = "@"
    $q4Sheet.Cells.Item($rowIdx, 2).Value2 = $rowData[1]
    $q4Sheet.Cells.Item($rowIdx, 3).Value2 = $rowData[2]
    $q4Sheet.Cells.Item($rowIdx, 4).Value2 = $rowData[3]
    $q4Sheet.Cells.Item($rowIdx, 5).Value2 = $rowData[4]
    $q4Sheet.Cells.Item($rowIdx, 6).Value2 = $rowData[5]
    $q4Sheet.Cells.Item($rowIdx, 7).Value2 = $rowData[6]

    $q4Sheet.Cells.Item($rowIdx, 8).Value2 = [int]$rowData[7]

    $rowIdx = $rowIdx + 1
}

# Step 3: update the "总计" (totals) sheet -- relabel the existing row as
# Q4 (with its new holding count) and append a new row carrying the Q3
# totals that used to live in that same row.
$totalSheet.Cells.Item(2, 2).Value2 = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value2 = 4
$totalSheet.Cells.Item(2, 4).Value2 = 0.19

$totalSheet.Cells.Item(2, 1).Copy()
$totalSheet.Cells.Item(3, 1).PasteSpecial(-4122)
$totalSheet.Cells.Item(3, 1).Value2 = 1
$totalSheet.Cells.Item(3, 2).Value2 = "2022-Q3"
$totalSheet.Cells.Item(3, 3).Value2 = 8
$totalSheet.Cells.Item(3, 4).Value2 = 0.19

# Restore the originally active tab (the "总计" totals sheet stayed the
# active tab across this edit).
$totalSheet.Activate()
